# Scheduled runner update: refresh computed profit figures on the
# per-job Leve sheets (currentAveragePrice* / LevePrice* / LeveProfit*
# columns, H:N) with newly pulled market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value2 = 5733.3335
$ws.Range("I112").Value2 = 686.6667
$ws.Range("J112").Value2 = 7175.2383
$ws.Range("K112").Value2 = 2060.0001
$ws.Range("L112").Value2 = 21525.7149
$ws.Range("M112").Value2 = -952.0001000000002
$ws.Range("N112").Value2 = -23741.7149
$ws.Range("H137").Value2 = 16062867
$ws.Range("I137").Value2 = 277594.56
$ws.Range("J137").Value2 = 50001200
$ws.Range("K137").Value2 = 832783.6799999999
$ws.Range("L137").Value2 = 150003600
$ws.Range("M137").Value2 = -830233.6799999999
$ws.Range("N137").Value2 = -150008700
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 763.74
$ws.Range("I32").Value2 = 769.6629
$ws.Range("J32").Value2 = 715.8182
$ws.Range("K32").Value2 = 769.6629
$ws.Range("L32").Value2 = 715.8182
$ws.Range("M32").Value2 = -482.6629
$ws.Range("N32").Value2 = -1289.8182
$ws.Range("H61").Value2 = 10104356
$ws.Range("I61").Value2 = 12349324
$ws.Range("J61").Value2 = 2000
$ws.Range("K61").Value2 = 12349324
$ws.Range("L61").Value2 = 2000
$ws.Range("M61").Value2 = -12349112
$ws.Range("N61").Value2 = -2424
$ws.Range("H74").Value2 = 675.7234
$ws.Range("I74").Value2 = 437
$ws.Range("J74").Value2 = 1457
$ws.Range("K74").Value2 = 437
$ws.Range("L74").Value2 = 1457
$ws.Range("M74").Value2 = 437
$ws.Range("N74").Value2 = -3205
$ws.Range("H77").Value2 = 675.7234
$ws.Range("I77").Value2 = 437
$ws.Range("J77").Value2 = 1457
$ws.Range("K77").Value2 = 2185
$ws.Range("L77").Value2 = 7285
$ws.Range("M77").Value2 = 2183
$ws.Range("N77").Value2 = -16021
$ws.Range("H136").Value2 = 10104356
$ws.Range("I136").Value2 = 12349324
$ws.Range("J136").Value2 = 2000
$ws.Range("K136").Value2 = 37047972
$ws.Range("L136").Value2 = 6000
$ws.Range("M136").Value2 = -37045422
$ws.Range("N136").Value2 = -11100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 8785957
$ws.Range("I134").Value2 = 10769506
$ws.Range("J134").Value2 = 1667.1428
$ws.Range("K134").Value2 = 32308518
$ws.Range("L134").Value2 = 5001.428400000001
$ws.Range("M134").Value2 = -32305983
$ws.Range("N134").Value2 = -10071.4284
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3792333
$ws.Range("I31").Value2 = 4274400
$ws.Range("J31").Value2 = 32211.4
$ws.Range("K31").Value2 = 4274400
$ws.Range("L31").Value2 = 32211.4
$ws.Range("M31").Value2 = -4274105
$ws.Range("N31").Value2 = -32801.4
$ws.Range("H34").Value2 = 3792333
$ws.Range("I34").Value2 = 4274400
$ws.Range("J34").Value2 = 32211.4
$ws.Range("K34").Value2 = 4274400
$ws.Range("L34").Value2 = 32211.4
$ws.Range("M34").Value2 = -4274198
$ws.Range("N34").Value2 = -32615.4
$ws.Range("H58").Value2 = 23463.285
$ws.Range("I58").Value2 = 27185.395
$ws.Range("J58").Value2 = 10605.091
$ws.Range("K58").Value2 = 27185.395
$ws.Range("L58").Value2 = 10605.091
$ws.Range("M58").Value2 = -26982.395
$ws.Range("N58").Value2 = -11011.091
$ws.Range("H132").Value2 = 3880.422
$ws.Range("I132").Value2 = 1644.8422
$ws.Range("J132").Value2 = 16016.429
$ws.Range("K132").Value2 = 4934.5266
$ws.Range("L132").Value2 = 48049.287
$ws.Range("M132").Value2 = -2404.5266
$ws.Range("N132").Value2 = -53109.287
$ws.Range("H134").Value2 = 24148904
$ws.Range("I134").Value2 = 28572342
$ws.Range("J134").Value2 = 6946649
$ws.Range("K134").Value2 = 85717026
$ws.Range("L134").Value2 = 20839947
$ws.Range("M134").Value2 = -85714491
$ws.Range("N134").Value2 = -20845017
$ws.Range("H136").Value2 = 23463.285
$ws.Range("I136").Value2 = 27185.395
$ws.Range("J136").Value2 = 10605.091
$ws.Range("K136").Value2 = 81556.185
$ws.Range("L136").Value2 = 31815.273
$ws.Range("M136").Value2 = -79006.185
$ws.Range("N136").Value2 = -36915.273
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 332.75
$ws.Range("I5").Value2 = 332.75
$ws.Range("K5").Value2 = 998.25
$ws.Range("M5").Value2 = -886.25
$ws.Range("H122").Value2 = 664.40625
$ws.Range("J122").Value2 = 999
$ws.Range("L122").Value2 = 8991
$ws.Range("N122").Value2 = -13891
$ws.Range("H129").Value2 = 1782.7142
$ws.Range("I129").Value2 = 460
$ws.Range("J129").Value2 = 2774.75
$ws.Range("K129").Value2 = 1380
$ws.Range("L129").Value2 = 8324.25
$ws.Range("M129").Value2 = 3620
$ws.Range("N129").Value2 = -18324.25
$ws.Range("H130").Value2 = 11206
$ws.Range("I130").Value2 = 50530
$ws.Range("J130").Value2 = 1375
$ws.Range("K130").Value2 = 151590
$ws.Range("L130").Value2 = 4125
$ws.Range("M130").Value2 = -146570
$ws.Range("N130").Value2 = -14165
$ws.Range("H131").Value2 = 12747058
$ws.Range("I131").Value2 = 58823904
$ws.Range("J131").Value2 = 1227846.5
$ws.Range("K131").Value2 = 176471712
$ws.Range("L131").Value2 = 3683539.5
$ws.Range("M131").Value2 = -176466672
$ws.Range("N131").Value2 = -3693619.5
$ws.Range("H135").Value2 = 332.75
$ws.Range("I135").Value2 = 332.75
$ws.Range("K135").Value2 = 2994.75
$ws.Range("M135").Value2 = -459.75
$ws.Range("H136").Value2 = 1543.3334
$ws.Range("I136").Value2 = 952
$ws.Range("J136").Value2 = 4500
$ws.Range("K136").Value2 = 2856
$ws.Range("L136").Value2 = 13500
$ws.Range("M136").Value2 = 2244
$ws.Range("N136").Value2 = -23700
$ws.Range("H139").Value2 = 31251704
$ws.Range("I139").Value2 = 50001030
$ws.Range("J139").Value2 = 2822.1667
$ws.Range("K139").Value2 = 150003090
$ws.Range("L139").Value2 = 8466.500100000001
$ws.Range("M139").Value2 = -149997950
$ws.Range("N139").Value2 = -18746.5001
$ws.Range("H140").Value2 = 1210.0358
$ws.Range("I140").Value2 = 973.087
$ws.Range("J140").Value2 = 2300
$ws.Range("K140").Value2 = 2919.261
$ws.Range("L140").Value2 = 6900
$ws.Range("M140").Value2 = 2260.739
$ws.Range("N140").Value2 = -17260
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 1241.6274
$ws.Range("I132").Value2 = 690.55554
$ws.Range("J132").Value2 = 1861.5834
$ws.Range("K132").Value2 = 2071.66662
$ws.Range("L132").Value2 = 5584.7502
$ws.Range("M132").Value2 = 458.33338
$ws.Range("N132").Value2 = -10644.7502
$ws.Range("H136").Value2 = 18871100
$ws.Range("I136").Value2 = 20836506
$ws.Range("J136").Value2 = 3201
$ws.Range("K136").Value2 = 62509518
$ws.Range("L136").Value2 = 9603
$ws.Range("M136").Value2 = -62506968
$ws.Range("N136").Value2 = -14703
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 19702406
$ws.Range("I132").Value2 = 20690384
$ws.Range("J132").Value2 = 17726448
$ws.Range("K132").Value2 = 62071152
$ws.Range("L132").Value2 = 53179344
$ws.Range("M132").Value2 = -62068622
$ws.Range("N132").Value2 = -53184404
$ws.Range("H136").Value2 = 26976564
$ws.Range("I136").Value2 = 22239862
$ws.Range("J136").Value2 = 62501830
$ws.Range("K136").Value2 = 66719586
$ws.Range("L136").Value2 = 187505490
$ws.Range("M136").Value2 = -66717036
$ws.Range("N136").Value2 = -187510590